$d = $word.ActiveDocument

# 1. Job title: "Entry Level Programmer" -> "QA Tester"
[void]$d.Content.Find.Execute("Entry Level Programmer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "QA Tester", 2)

# 2 & 4. Company name: "IO Interactive" -> "Ubisoft" (both occurrences)
[void]$d.Content.Find.Execute("IO Interactive", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ubisoft", 2)

# 3. Replace the "looking forward" sentence about Hitman/Project 007 with the Far Cry paragraph
$oldSentence = "I am really looking forward to Hitman 3 next year and gameplay footage of Project 007, both as a gamer and game developer."
$newSentence = "I have played all the entries in the Far Cry series and I have been a fan since the first Far Cry on CD. I was really surprised by the two sequels that followed it as they are all completely unique and standout from each other. My favourite entries are Far Cry 2 and 3. Vaas was such a great antagonist. I am also really looking forward to Beyond Good and Evil 2."
[void]$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newSentence, 2)
